$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update company name for REC Silicon ASA in row 3
$ws.Range("B3").Value = "REC Silicon ASA (OB:RECSI)"

# Values that are identical across rows 2 and 3
$updates = @{
    "D" = -0.209
    "G" = -0.1375212224108659
    "H" = -0.1528013582342954
    "I" = -0.3516079613987593
    "J" = -0.3516079613987593
    "K" = -65.59999999999999
    "L" = -0.5568760611205432
    "U" = 35.9
    "V" = 0.05132237312365975
    "W" = -1.547169811320755
    "X" = 0.1110990850464327
    "Y" = -1.658268896367187
    "Z" = 0.6651718584966781
    "AA" = -0.233879721145841
    "AB" = 0.09008362419156583
    "AC" = -0.3239633453374068
    "AD" = 220.3
    "AE" = 0.09708926386920136
    "AF" = 220.3970892638692
    "AG" = 184.4970892638692
    "AH" = 0.239588853836072
    "AI" = 1.095428815944843
    "AJ" = 0.2087078017615482
    "AK" = 1.116154495433071
    "AL" = 19.9
    "AM" = 19.5
    "AN" = -314.7142857142857
    "AO" = -2.085427135678392
    "AP" = -263.567270376956
    "AQ" = -2.128205128205128
}

foreach ($col in $updates.Keys) {
    $value = $updates[$col]
    $ws.Range("${col}2").Value = $value
    $ws.Range("${col}3").Value = $value
}
